$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4556.9067
$ws.Range("I40").Value = 2711.8
$ws.Range("J40").Value = 4799.684
$ws.Range("K40").Value = 2711.8
$ws.Range("L40").Value = 4799.684
$ws.Range("M40").Value = -2536.8
$ws.Range("N40").Value = -5149.684

$ws.Range("H64").Value = 5232.5
$ws.Range("I64").Value = 4709.2856
$ws.Range("J64").Value = 5755.7144
$ws.Range("K64").Value = 4709.2856
$ws.Range("L64").Value = 5755.7144
$ws.Range("M64").Value = -4461.2856
$ws.Range("N64").Value = -6251.7144

$ws.Range("H67").Value = 5232.5
$ws.Range("I67").Value = 4709.2856
$ws.Range("J67").Value = 5755.7144
$ws.Range("K67").Value = 4709.2856
$ws.Range("L67").Value = 5755.7144
$ws.Range("M67").Value = -3851.2856
$ws.Range("N67").Value = -7471.7144

$ws.Range("H113").Value = 3196.5417
$ws.Range("J113").Value = 3605.875
$ws.Range("L113").Value = 3605.875
$ws.Range("N113").Value = -10113.875

$ws.Range("H132").Value = 17674.666
$ws.Range("I132").Value = 17674.666
$ws.Range("K132").Value = 53023.99800000001
$ws.Range("M132").Value = -50493.99800000001

$ws.Range("H137").Value = 18525308
$ws.Range("I137").Value = 33335006
$ws.Range("J137").Value = 13182.667
$ws.Range("K137").Value = 100005018
$ws.Range("L137").Value = 39548.001
$ws.Range("M137").Value = -100002468
$ws.Range("N137").Value = -44648.001


# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3712689.2
$ws.Range("I61").Value = 13957.444
$ws.Range("K61").Value = 13957.444
$ws.Range("M61").Value = -13745.444

$ws.Range("H122").Value = 1593.2222
$ws.Range("I122").Value = 1475.1765
$ws.Range("J122").Value = 3600
$ws.Range("K122").Value = 4425.529500000001
$ws.Range("L122").Value = 10800
$ws.Range("M122").Value = -1975.529500000001
$ws.Range("N122").Value = -15700

$ws.Range("H135").Value = 100000.336
$ws.Range("J135").Value = 100000.336
$ws.Range("L135").Value = 100000.336
$ws.Range("N135").Value = -110140.336

$ws.Range("H136").Value = 3712689.2
$ws.Range("I136").Value = 13957.444
$ws.Range("K136").Value = 41872.33199999999
$ws.Range("M136").Value = -39322.33199999999


# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 68204.164
$ws.Range("J135").Value = 68204.164
$ws.Range("L135").Value = 68204.164
$ws.Range("N135").Value = -78344.164


# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5264635
$ws.Range("J31").Value = 562.3333
$ws.Range("L31").Value = 562.3333
$ws.Range("N31").Value = -1152.3333

$ws.Range("H34").Value = 5264635
$ws.Range("J34").Value = 562.3333
$ws.Range("L34").Value = 562.3333
$ws.Range("N34").Value = -966.3333

$ws.Range("H107").Value = 641.1739
$ws.Range("I107").Value = 462.4
$ws.Range("K107").Value = 462.4
$ws.Range("M107").Value = 1457.6


# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 130.28572
$ws.Range("J38").Value = 208.42857
$ws.Range("L38").Value = 625.28571
$ws.Range("N38").Value = -1319.28571

$ws.Range("H107").Value = 4074.9285
$ws.Range("J107").Value = 5065.2
$ws.Range("L107").Value = 15195.6
$ws.Range("N107").Value = -19035.6

$ws.Range("I126").Value = 15000
$ws.Range("K126").Value = 45000
$ws.Range("M126").Value = -40060

$ws.Range("H129").Value = 528700.7
$ws.Range("I129").Value = 1001734.3
$ws.Range("J129").Value = 3107.7778
$ws.Range("K129").Value = 3005202.9
$ws.Range("L129").Value = 9323.3334
$ws.Range("M129").Value = -3000202.9
$ws.Range("N129").Value = -19323.3334

$ws.Range("H137").Value = 9413.333000000001
$ws.Range("I137").Value = 3107.25
$ws.Range("K137").Value = 9321.75
$ws.Range("M137").Value = -4221.75


# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 57250
$ws.Range("J134").Value = 57250
$ws.Range("L134").Value = 171750
$ws.Range("N134").Value = -176820


# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3602.875
$ws.Range("I7").Value = 3388.0908
$ws.Range("J7").Value = 4075.4
$ws.Range("K7").Value = 3388.0908
$ws.Range("L7").Value = 4075.4
$ws.Range("M7").Value = -3276.0908
$ws.Range("N7").Value = -4299.4

$ws.Range("H22").Value = 2811.625
$ws.Range("I22").Value = 2010.625
$ws.Range("J22").Value = 3212.125
$ws.Range("K22").Value = 2010.625
$ws.Range("L22").Value = 3212.125
$ws.Range("M22").Value = -1715.625
$ws.Range("N22").Value = -3802.125

$ws.Range("H27").Value = 2811.625
$ws.Range("I27").Value = 2010.625
$ws.Range("J27").Value = 3212.125
$ws.Range("K27").Value = 2010.625
$ws.Range("L27").Value = 3212.125
$ws.Range("M27").Value = -1903.625
$ws.Range("N27").Value = -3426.125

$ws.Range("H40").Value = 3853.238
$ws.Range("I40").Value = 2038.6428
$ws.Range("K40").Value = 2038.6428
$ws.Range("M40").Value = -1902.6428

$ws.Range("H93").Value = 5043.6
$ws.Range("I93").Value = 3917.4285
$ws.Range("K93").Value = 3917.4285
$ws.Range("M93").Value = -2669.4285

$ws.Range("H122").Value = 4626.5
$ws.Range("I122").Value = 2956.2144
$ws.Range("J122").Value = 6964.9
$ws.Range("K122").Value = 8868.643199999999
$ws.Range("L122").Value = 20894.7
$ws.Range("M122").Value = -6418.643199999999
$ws.Range("N122").Value = -25794.7

$ws.Range("H126").Value = 3602.875
$ws.Range("I126").Value = 3388.0908
$ws.Range("J126").Value = 4075.4
$ws.Range("K126").Value = 10164.2724
$ws.Range("L126").Value = 12226.2
$ws.Range("M126").Value = -7694.2724
$ws.Range("N126").Value = -17166.2

$ws.Range("H138").Value = 29990
$ws.Range("J138").Value = 29990
$ws.Range("L138").Value = 29990
$ws.Range("M138").Value = -40270

$ws.Range("H140").Value = 100429
$ws.Range("J140").Value = 100429
$ws.Range("L140").Value = 100429
$ws.Range("M140").Value = -110789


# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 8515.412
$ws.Range("I81").Value = 1994.1818
$ws.Range("J81").Value = 20471
$ws.Range("K81").Value = 3988.3636
$ws.Range("L81").Value = 40942
$ws.Range("M81").Value = -2927.3636
$ws.Range("N81").Value = -43064

$ws.Range("H84").Value = 8515.412
$ws.Range("I84").Value = 1994.1818
$ws.Range("J84").Value = 20471
$ws.Range("K84").Value = 19941.818
$ws.Range("L84").Value = 204710
$ws.Range("M84").Value = -14637.818
$ws.Range("N84").Value = -215318

$ws.Range("H126").Value = 1658.5264
$ws.Range("J126").Value = 2005.375
$ws.Range("L126").Value = 6016.125
$ws.Range("N126").Value = -10956.125

